$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-completed practice programs (rows 4-7, columns B & C)
$ws.Range("B4").Value = "multiply2Floating.java"
$ws.Range("C4").Value = "DONE"

$ws.Range("B5").Value = "findASCII.java"
$ws.Range("C5").Value = "DONE"

$ws.Range("B6").Value = "findQuotientAndRemainder.java"
$ws.Range("C6").Value = "DONE"

$ws.Range("B7").Value = "swapUsing3Variable.java"
$ws.Range("C7").Value = "DONE"

# Update the view: scroll so row 4 is the top-left visible row, and select C7
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C7").Select()
